# This script applies the "correção nos dados e inicio da analise PNAD 2009"
# edit: the worksheet contained category header rows ("sexo", "cor ou raça",
# "grupos de idade", "classes de rendimento mensal domiciliar per capita")
# and a trailing source-note row ("fonte: ...") that had no numeric data of
# their own. The fix removes those five empty/placeholder rows entirely so
# that the numeric data that follows each one shifts up into a contiguous
# table. Deleting the rows (instead of just clearing cells) also prunes the
# now-unused shared strings automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so the row numbers of rows still to be removed
# do not shift as we go.
$ws.Rows("26").Delete()   # "fonte: ibge, diretoria de pesquisas, ..."
$ws.Rows("19").Delete()   # "classes de rendimento mensal domiciliar per capita"
$ws.Rows("13").Delete()   # "grupos de idade"
$ws.Rows("8").Delete()    # "cor ou raça"
$ws.Rows("5").Delete()    # "sexo"
